# Update countries & provincias Spain
# Applies the data refresh captured in the commit "Update countries & provincias Spain":
#  - Re-ranked a few countries (Republica de Chipre overtakes Liberia/Gambia;
#    Timor Oriental overtakes Santa Lucia) causing their table rows to swap places.
#  - Refreshed totals (Casos totales, Nuevos casos, Casos activos, Recuperados,
#    Casos criticos, Muertes) for the affected rows.
#  - Bumped the "Datos actualizados" timestamp footer.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Footer timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 10 de Agosto de 2020 a las 02:18"

# --- Helper: write a full data row (B:H) -------------------------------
function Set-Row($row, $b, $c, $d, $e, $f, $g, $h) {
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 6).Value = $f
    $ws.Cells.Item($row, 7).Value = $g
    $ws.Cells.Item($row, 8).Value = $h
}

# Row 4 - Estados Unidos
Set-Row 4 5199411 47816 2664612 2369191 0 525 165608

# Row 5 - Brasil
Set-Row 5 3035582 22213 2118460 815986 0 593 101136

# Row 20 - Argentina
Set-Row 20 246499 4688 108242 133651 0 83 4606

# Row 27 - Canada
Set-Row 27 119451 230 103728 6742 0 5 8981

# Row 39 - Panama
Set-Row 39 74492 841 48748 24105 0 30 1639

# Row 83 - Sudan
Set-Row 83 11956 62 6266 4909 0 8 781

# Row 96 - Luxemburgo
Set-Row 96 7205 36 5848 1237 0 0 120

# Row 110 - Republica de Africa Central
Set-Row 110 4641 0 1721 2860 0 1 60

# Row 114 - Montenegro
Set-Row 114 3618 30 2452 1102 0 2 64

# Row 126 - Surinam
Set-Row 126 2391 85 1635 727 0 0 29

# Rows 144-146: Republica de Chipre overtakes Liberia and Gambia.
# Row 144 becomes Republica de Chipre with fresh totals; Liberia and Gambia
# shift down one row each, keeping their previous (unrefreshed) totals.
$ws.Cells.Item(144, 1).Value = "Republica de Chipre"
Set-Row 144 1242 9 870 353 0 0 19

$ws.Cells.Item(145, 1).Value = "Liberia"
Set-Row 145 1237 3 723 435 0 0 79

$ws.Cells.Item(146, 1).Value = "Gambia"
Set-Row 146 1235 145 221 991 0 4 23

# Row 163 - Reunion
Set-Row 163 687 6 631 51 0 0 5

# Rows 202-203: Timor Oriental overtakes Santa Lucia (same totals, so only
# the country names swap).
$ws.Cells.Item(202, 1).Value = "Timor Oriental"
$ws.Cells.Item(203, 1).Value = "Santa Lucia"
